$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "123"
$ws.Range("C7").Value = 0

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "213"
$ws.Range("C8").Value = 0

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = 0

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "2"
$ws.Range("C10").Value = 0
